$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "general": update a handful of summary values
# ---------------------------------------------------------------------------
$wsGeneral = $wb.Worksheets.Item("general")
$wsGeneral.Range("B3").Value = 30.19742239378569
$wsGeneral.Range("B4").Value = 0.01200008392333984
$wsGeneral.Range("B6").Value = 30.16742239378572
$wsGeneral.Range("B7").Value = 0
$wsGeneral.Range("B8").Value = 0
$wsGeneral.Range("B9").Value = 0.02999999999997272
$wsGeneral.Range("B10").Value = 0

# ---------------------------------------------------------------------------
# Sheet "x": update values (no structural change)
# ---------------------------------------------------------------------------
$wsX = $wb.Worksheets.Item("x")
$wsX.Range("B2").Value = 4
$wsX.Range("B3").Value = 6
$wsX.Range("B4").Value = 10
$wsX.Range("B7").Value = 13
$wsX.Range("B9").Value = 12
$wsX.Range("B10").Value = 3
$wsX.Range("B13").Value = 5

# ---------------------------------------------------------------------------
# Sheet "U": update values (no structural change)
# ---------------------------------------------------------------------------
$wsU = $wb.Worksheets.Item("U")
$wsU.Range("B4").Value = 3
$wsU.Range("B6").Value = 3
$wsU.Range("B7").Value = 3
$wsU.Range("B13").Value = 3

# ---------------------------------------------------------------------------
# Sheet "TBar": update values (no structural change)
# ---------------------------------------------------------------------------
$wsTBar = $wb.Worksheets.Item("TBar")
$wsTBar.Range("B5").Value = 20.34885527085025
$wsTBar.Range("B7").Value = 26.71579249669673
$wsTBar.Range("B8").Value = 20
$wsTBar.Range("B13").Value = 29.31314932600328
$wsTBar.Range("B14").Value = 27.90090852477161
$wsTBar.Range("B15").Value = 22.66758337047728

# ---------------------------------------------------------------------------
# Sheet "y": remove all data rows except the header (rows 2:10 deleted)
# ---------------------------------------------------------------------------
$wsY = $wb.Worksheets.Item("y")
$wsY.Rows("2:10").Delete()

# ---------------------------------------------------------------------------
# Sheet "Q": update every data value (rows 7:71)
# ---------------------------------------------------------------------------
$wsQ = $wb.Worksheets.Item("Q")
$wsQ.Range("C7").Value = 250.9700000000009
$wsQ.Range("C8").Value = 260.9900000000009
$wsQ.Range("C9").Value = 252.9750000000009
$wsQ.Range("C10").Value = 269.5799999999992
$wsQ.Range("C11").Value = 250.5750000000009
$wsQ.Range("C12").Value = 332.4450000000016
$wsQ.Range("C13").Value = 341.6700000000017
$wsQ.Range("C14").Value = 337.6900000000016
$wsQ.Range("C15").Value = 350.015
$wsQ.Range("C16").Value = 337.9200000000017
$wsQ.Range("C17").Value = 154.3
$wsQ.Range("C18").Value = 148.3449999999993
$wsQ.Range("C19").Value = 128.7049999999993
$wsQ.Range("C20").Value = 146.3249999999993
$wsQ.Range("C21").Value = 134.2149999999993
$wsQ.Range("C22").Value = 72.6299999999995
$wsQ.Range("C23").Value = 80.0549999999995
$wsQ.Range("C24").Value = 82.31999999999948
$wsQ.Range("C25").Value = 83.9549999999995
$wsQ.Range("C26").Value = 80.8149999999995
$wsQ.Range("C27").Value = 295.9199999999996
$wsQ.Range("C28").Value = 323.5
$wsQ.Range("C29").Value = 294.2649999999996
$wsQ.Range("C30").Value = 311.1
$wsQ.Range("C31").Value = 297.3649999999997
$wsQ.Range("C32").Value = 107.3799999999999
$wsQ.Range("C33").Value = 112.2399999999999
$wsQ.Range("C34").Value = 93.78999999999985
$wsQ.Range("C35").Value = 108.8349999999998
$wsQ.Range("C36").Value = 94.77999999999986
$wsQ.Range("C37").Value = 141.0250000000001
$wsQ.Range("C38").Value = 143.4
$wsQ.Range("C39").Value = 139.7050000000001
$wsQ.Range("C40").Value = 150.4250000000002
$wsQ.Range("C41").Value = 134.7700000000002
$wsQ.Range("C42").Value = 140.5549999999989
$wsQ.Range("C43").Value = 159.214999999999
$wsQ.Range("C44").Value = 142.1399999999989
$wsQ.Range("C45").Value = 147.7249999999989
$wsQ.Range("C46").Value = 139.7449999999989
$wsQ.Range("C47").Value = 226.0399999999994
$wsQ.Range("C48").Value = 247.1799999999994
$wsQ.Range("C49").Value = 221.8549999999994
$wsQ.Range("C50").Value = 238.4549999999994
$wsQ.Range("C51").Value = 224.4749999999994
$wsQ.Range("C52").Value = 57.95
$wsQ.Range("C53").Value = 58.67999999999927
$wsQ.Range("C54").Value = 61.72999999999927
$wsQ.Range("C55").Value = 60.65499999999928
$wsQ.Range("C56").Value = 52.91499999999927
$wsQ.Range("C57").Value = 332.4450000000016
$wsQ.Range("C58").Value = 341.6700000000017
$wsQ.Range("C59").Value = 337.6900000000016
$wsQ.Range("C60").Value = 350.015
$wsQ.Range("C61").Value = 337.9200000000017
$wsQ.Range("C62").Value = 295.9199999999996
$wsQ.Range("C63").Value = 323.5
$wsQ.Range("C64").Value = 294.2649999999996
$wsQ.Range("C65").Value = 311.1
$wsQ.Range("C66").Value = 297.3649999999997
$wsQ.Range("C67").Value = 154.3
$wsQ.Range("C68").Value = 148.3449999999993
$wsQ.Range("C69").Value = 128.7049999999993
$wsQ.Range("C70").Value = 146.3249999999993
$wsQ.Range("C71").Value = 134.2149999999993

# ---------------------------------------------------------------------------
# Sheet "R": zero out the remaining (nonzero) entries
# ---------------------------------------------------------------------------
$wsR = $wb.Worksheets.Item("R")
$wsR.Range("C2").Value = 0
$wsR.Range("C3").Value = 0
$wsR.Range("C4").Value = 0
$wsR.Range("C5").Value = 0
$wsR.Range("C6").Value = 0
$wsR.Range("C12").Value = 0
$wsR.Range("C13").Value = 0
$wsR.Range("C15").Value = 0
$wsR.Range("C16").Value = 0

# ---------------------------------------------------------------------------
# Sheet "L": zero out the remaining (nonzero) entries
# ---------------------------------------------------------------------------
$wsL = $wb.Worksheets.Item("L")
$wsL.Range("C12").Value = 0
$wsL.Range("C13").Value = 0
$wsL.Range("C14").Value = 0
$wsL.Range("C15").Value = 0
$wsL.Range("C16").Value = 0
$wsL.Range("C22").Value = 0
$wsL.Range("C23").Value = 0
$wsL.Range("C24").Value = 0
$wsL.Range("C25").Value = 0
$wsL.Range("C26").Value = 0
$wsL.Range("C27").Value = 0
$wsL.Range("C28").Value = 0
$wsL.Range("C29").Value = 0
$wsL.Range("C30").Value = 0
$wsL.Range("C31").Value = 0

# ---------------------------------------------------------------------------
# Sheet "rho": keep row 2 (updated) but drop rows 3:10
# ---------------------------------------------------------------------------
$wsRho = $wb.Worksheets.Item("rho")
$wsRho.Rows("3:10").Delete()
$wsRho.Range("B2").Value = 4

# ---------------------------------------------------------------------------
# Sheet "alpha": remove all data rows except the header (rows 2:10 deleted)
# ---------------------------------------------------------------------------
$wsAlpha = $wb.Worksheets.Item("alpha")
$wsAlpha.Rows("2:10").Delete()
